# Update the MSME Tajikistan Summary sheet with refreshed (more precise)
# figures for "Enterprises density (per 1000 people)" and
# "Enterprises (% of total)".
#
# A leading apostrophe is used so Excel stores these numeric-looking values
# as literal text (matching the workbook's existing shared-string / text
# cell type for these rows) instead of silently converting them to numbers.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Enterprises density (per 1000 people) -- row 13
$ws.Range("B13").Value = "'20.72"
$ws.Range("C13").Value = "'1.12"
$ws.Range("D13").Value = "'21.84"

# Enterprises (% of total) -- row 16
$ws.Range("B16").Value = "'94.79"
$ws.Range("C16").Value = "'5.14"
$ws.Range("D16").Value = "'99.92"
